$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 301.5198272274479
$ws.Range("G2").Value = 20.9257198412324
$ws.Range("H2").Value = 576.1834978760137
$ws.Range("I2").Value = 0.2408741705832104
$ws.Range("J2").Value = 0.001367532199953357
$ws.Range("K2").Value = 0.5638990011166842
$ws.Range("L2").Value = 0.1807565214380907
$ws.Range("M2").Value = 0.00810850077414255
$ws.Range("N2").Value = 0.3673345428161195

$ws.Range("F3").Value = 0.002623255133635339
$ws.Range("G3").Value = 0.00190431164200106
$ws.Range("H3").Value = 0.003374091708301619
$ws.Range("I3").Value = 0.002421097065732277
$ws.Range("J3").Value = 0.001751948135999053
$ws.Range("K3").Value = 0.003115281650776451
$ws.Range("L3").Value = 0.002749737008511098
$ws.Range("M3").Value = 0.002013126224807544
$ws.Range("N3").Value = 0.003517830068651254

$ws.Range("F4").Value = 301.5224504825815
$ws.Range("G4").Value = 20.9276241528744
$ws.Range("H4").Value = 576.1868719677219
$ws.Range("I4").Value = 0.2432952676489426
$ws.Range("J4").Value = 0.00311948033595241
$ws.Range("K4").Value = 0.5670142827674607
$ws.Range("L4").Value = 0.1835062584466018
$ws.Range("M4").Value = 0.01012162699895009
$ws.Range("N4").Value = 0.3708523728847709
